# Scale the parse to 20 objects: append two more "category" rows
# (GNU and Burton brands) to the Categories sheet, directly below
# the existing K2 row (row 21), following the same column layout:
# data-object | key | description.en-US | externalId | name.en-US |
# slug.en-US | parent.key | parent.typeId
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "category"
$ws.Range("B22").Value = "GNUKey"
$ws.Range("C22").Value = "GNUDescription"
$ws.Range("D22").Value = "GNUId"
$ws.Range("E22").Value = "GNU"
$ws.Range("F22").Value = "GNUSlug"
$ws.Range("G22").Value = "brandKey"
$ws.Range("H22").Value = "category"

$ws.Range("A23").Value = "category"
$ws.Range("B23").Value = "BurtonKey"
$ws.Range("C23").Value = "BurtonDescription"
$ws.Range("D23").Value = "BurtonId"
$ws.Range("E23").Value = "Burton"
$ws.Range("F23").Value = "BurtonSlug"
$ws.Range("G23").Value = "brandKey"
$ws.Range("H23").Value = "category"
